# Split the "Game logic is moved to separated classes." bullet from the
# stray "_GoBack" bookmark that currently sits between its two runs, and
# add a new bullet "Removed old files that are unused." right after it,
# with the "_GoBack" bookmark now collapsed at the end of the new bullet.

$d = $word.ActiveDocument

# 1) Drop the existing "_GoBack" bookmark; this leaves "Game l" and
#    "ogic is moved to separated classes." as plain adjacent runs.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2) Find the (now last) paragraph — "Game logic is moved to separated
#    classes." — and add a fresh bullet paragraph right after it,
#    inheriting the same ListParagraph / numbering formatting.
$gameLogicPara = $d.Paragraphs.Last
$gameLogicPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last

# 3) Fill the new paragraph with its text. A trailing placeholder
#    character is appended temporarily so the bookmark can be anchored
#    immediately after the real text without landing on the paragraph's
#    end-of-range (pilcrow) slot.
$newPara.Range.InsertAfter("Removed old files that are unused.X")

# 4) Re-create "_GoBack", collapsed, right before the placeholder — i.e.
#    right after the real text.
$markerPos = $newPara.Range.End - 2
$markerRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $markerRange)

# 5) Remove the placeholder character, leaving the bookmark collapsed at
#    the tail of the paragraph's text.
$placeholderPos = $newPara.Range.End - 2
$placeholderRange = $d.Range($placeholderPos, $placeholderPos + 1)
$placeholderRange.Delete()
